$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.145.41'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.80%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.341.85'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +0.44%  '

# Row 4
$ws.Range('E4').Value = '  +0.11%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '583.75'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.73%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '177.09'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.91%  '

# Row 7
$ws.Range('E7').Value = '  +0.00%  '

# Row 8
$ws.Range('E8').Value = '  +0.61%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.182'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +4.22%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.582'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +1.26%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '47.99'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +5.84%  '

# Row 12
$ws.Range('E12').Value = '  +1.46%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '692.85'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +3.84%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.887.08'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.76%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.41'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.33%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '68.223.81'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.74%  '

# Row 17
$ws.Range('E17').Value = '  +1.26%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.343.32'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.76%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.41'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.07%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.17'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +2.54%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.894'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.79%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.45'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.31%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.01'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.40%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '100.20'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.89%  '

# Row 25
$ws.Range('E25').Value = '  +2.43%  '

# Row 26
$ws.Range('E26').Value = '  +1.16%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.52'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +2.54%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '32.95'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -2.67%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.50'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +1.18%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.94'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -5.62%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '565.39'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -3.89%  '

# Row 32
$ws.Range('E32').Value = '  +1.11%  '

# Row 33
$ws.Range('E33').Value = '  +1.45%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '57.43'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +1.20%  '

# Row 35
$ws.Range('E35').Value = '  -0.03%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.691.28'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.19%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.26'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.13%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.136'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +4.29%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '34.73'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +5.62%  '

# Row 40
$ws.Range('E40').Value = '  +2.62%  '

# Row 41
$ws.Range('E41').Value = '  -0.43%  '

# Row 42
$ws.Range('B42').Value = 'TheGraph'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.335'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.88%  '

# Row 43
$ws.Range('B43').Value = 'PEPE'
$ws.Range('C43').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0₃0671'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.97%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.28'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +1.03%  '

# Row 45
$ws.Range('E45').Value = '  +2.16%  '

# Row 46
$ws.Range('E46').Value = '  +2.75%  '

# Row 47
$ws.Range('E47').Value = '  +0.79%  '

# Row 48
$ws.Range('E48').Value = '  +0.07%  '

# Row 49
$ws.Range('E49').Value = '  -0.06%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '130.90'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +3.24%  '

# Row 51
$ws.Range('E51').Value = '  +1.04%  '
